# Update the "想去人数" (want-to-go headcount) column (F) figures on the
# "展览" and "全部类型" sheets to the freshly scraped totals.
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3091
    5  = 2556
    7  = 130
    9  = 1304
    11 = 54
    13 = 1153
    14 = 331
    15 = 321
    18 = 105
    20 = 86
    21 = 2332
    22 = 20
    23 = 270
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
